$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H; this shifts the existing H:L
# (PerturbationDE..NumProcessors) one column to the right, to I:M.
$ws.Columns("H:H").Insert()

# Give the freshly inserted column roughly the same width as its
# neighbour (column G) -- closest value this host's pixel-quantised
# ColumnWidth setter can reach.
$ws.Columns("H").ColumnWidth = 9.5

# New header + value for the "ParetoSize" parameter living in the
# column that was just inserted.
$ws.Range("H1").Value = "ParetoSize"
$ws.Range("H2").Value = 10

# A couple of other parameter values changed at the same time.
$ws.Range("B2").Value = 2
$ws.Range("G2").Value = 10

# The user's cell selection ended up on H3 after the edit.
[void]$ws.Range("H3").Select()
